$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows ---
# C5: 1.04 -> 2.6
$ws.Range("C5").Value = 2.6

# Extend the SUM formula in E2 to cover more rows
$ws.Range("E2").Formula = "=SUM(C2:C20)"

# --- Fill in the previously-empty rows 8-10 (already had currency style) ---
$ws.Range("A8").Value = "PCB "
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 5

$ws.Range("A9").Value = "header 8 pinos"
$ws.Range("B9").Value = 3
$ws.Range("C9").Value = 0.3

$ws.Range("A10").Value = "header 6 pinos"
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 0.2

# --- New rows 11-12: apply the same currency number format used by C2:C10 ---
$currencyFormat = "#,##0.00\ ""€"""

$ws.Range("A11").Value = "conector jst"
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = 0.85
$ws.Range("C11").NumberFormat = $currencyFormat

$ws.Range("A12").Value = "bateria"
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = 6
$ws.Range("C12").NumberFormat = $currencyFormat

# --- New rows 13-14: left as plain numbers (no currency format), matching source ---
$ws.Range("A13").Value = "carregador da bateria"
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = 8.5

$ws.Range("A14").Value = "pente 40 pinos"
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = 0.65

# --- Selection moves to E3 ---
$ws.Range("E3").Select()
